$wb = $excel.ActiveWorkbook

# Add the new "Tooltip" worksheet after the last existing sheet (RoundTripFields)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Tooltip"

# Row 1: bold header "ToolTip Text"
$newSheet.Range("A1").Value = "ToolTip Text"
$newSheet.Range("A1").Font.Bold = $true

# Row 2: wrapped long description text
$newSheet.Range("A2").Value = "An engagement is typically considered a potential round trip if it is acquired by a sponsor (subject is a potential round trip) or by a sponsor-backed operating company (buyer is a potential round trip). Note ""sponsor"" includes firms tagged as Private Equity Group, Hedge Fund, or Family Office."
$newSheet.Range("A2").WrapText = $true
$newSheet.Rows.Item(2).RowHeight = 72

# Column width for A
$newSheet.Columns.Item(1).ColumnWidth = 56.83

# Match the selection / active cell state recorded in the target workbook
[void]$newSheet.Range("A7").Select()
